$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly report data between row 2 and row 3 for the
# date/volume/price columns (D, M, N, O, P, S).

# Row 2 new values (previously row 3's values)
$ws.Range("D2").Value = 44421
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("S2").Value = 1200

# Row 3 new values (previously row 2's values)
$ws.Range("D3").Value = 44291
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 23000
$ws.Range("P3").Value = 23000
$ws.Range("S3").Value = 1150
